$wb = $excel.ActiveWorkbook

# --- Build the new "agent2sensors" sheet by copying the existing one ------
# (copying -- rather than Worksheets.Add -- keeps the same sheetFormatPr /
# default row height as the rest of the workbook). The copy is placed right
# before "targets", i.e. it becomes the first sheet in the workbook.
$targetsWs = $wb.Worksheets.Item("targets")
$sensorsWs = $wb.Worksheets.Item("agent2sensors")
$sensorsWs.Copy($targetsWs) | Out-Null

# NOTE: worksheet variables captured before a sheet is inserted/removed can
# end up pointing at the wrong tab afterwards, so re-resolve everything by
# name once the sheet collection has changed shape.
$oldSensorsWs = $wb.Worksheets.Item("agent2sensors")
$newSensorsWs = $wb.Worksheets.Item("agent2sensors (2)")

# Rename the old sheet out of the way, then claim its old name for the copy.
$oldSensorsWs.Name = "agent2sensorsOld"
$newSensorsWs.Name = "agent2sensors"

# Re-resolve again now that both names changed.
$oldSensorsWs = $wb.Worksheets.Item("agent2sensorsOld")
$newSensorsWs = $wb.Worksheets.Item("agent2sensors")
$targetsWs = $wb.Worksheets.Item("targets")

# --- Replace the copied sheet's data with the new, smaller, generic table -
$newSensorsWs.Cells.Clear() | Out-Null

$headers = @("sensor1","sensor2","sensor3","sensor4","sensor5","sensor6","sensor7","sensor8")
for ($col = 1; $col -le 8; $col++) {
    $newSensorsWs.Cells.Item(1, $col).Value = $headers[$col - 1]
}

for ($row = 2; $row -le 6; $row++) {
    for ($col = 1; $col -le 8; $col++) {
        $newSensorsWs.Cells.Item($row, $col).Value = 1
    }
}

# Page setup to match a freshly printed sheet (A4 / portrait).
$newSensorsWs.PageSetup.PaperSize = 9
$newSensorsWs.PageSetup.Orientation = 1

# Selection: H2:H6 with the active cell at the top (H2).
$newSensorsWs.Activate()
$newSensorsWs.Range("H2:H6").Select() | Out-Null

# --- Update the selections left behind on the other two sheets ------------
$targetsWs.Activate()
$targetsWs.Range("A2:E30").Select() | Out-Null

$oldSensorsWs.Activate()
$oldSensorsWs.Range("J18").Select() | Out-Null

# --- Leave "agent2sensors" as the active sheet/tab ------------------------
$newSensorsWs.Activate()
